$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching source formatting)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '54.242.39'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '2.278.40'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').Value = '497.53'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = '128.46'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.0952'
$ws.Range('E9').Value = '  +2.20%  '
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('D12').Value = '4.71'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('D13').Value = '2.683.06'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '22.55'
$ws.Range('E14').Value = '  +4.93%  '
$ws.Range('D15').Value = '54.192.14'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '2.275.79'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('E18').Value = '  +4.16%  '
$ws.Range('D19').Value = '4.11'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').Value = '303.72'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '61.78'
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +2.01%  '
$ws.Range('E26').Value = '  +2.22%  '
$ws.Range('D27').Value = '174.98'
$ws.Range('E27').Value = '  +7.26%  '
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = '0.924'
$ws.Range('E35').Value = '  +9.21%  '
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('E37').Value = '  +2.67%  '
$ws.Range('D38').Value = '0.374'
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('E40').Value = '  +1.52%  '
$ws.Range('D41').Value = '124.97'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('D42').Value = '4.76'
$ws.Range('E42').Value = '  -2.43%  '
$ws.Range('E43').Value = '  +2.74%  '
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Value = '240.01'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('E49').Value = '  +1.07%  '
$ws.Range('D50').Value = '16.25'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('E51').Value = '  +0.30%  '
